$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: account holder first/last name swap
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must stay stored as TEXT (like the
# original inline string), not be auto-coerced into a numeric value. Force
# text with a leading apostrophe, then restore the original (non quote-
# prefixed) cell formatting by pasting formats in from a sibling cell that
# already carries the correct style.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("B6").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 09.10.2024"

# Row 6 (existing transaction, shifted values)
$ws.Range("B6").Value = "13.10."
$ws.Range("C6").Value = "14.10."
$ws.Range("D6").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E6").Value = "66,63-"

# Row 7
$ws.Range("B7").Value = "16.10."
$ws.Range("C7").Value = "17.10."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-53103382"
$ws.Range("E7").Value = "54,52-"

# Row 8
$ws.Range("B8").Value = "20.10."
$ws.Range("C8").Value = "21.10."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "25,33-"

# Row 9
$ws.Range("B9").Value = "23.10."
$ws.Range("C9").Value = "24.10."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 75735962"
$ws.Range("E9").Value = "84,38-"

# Row 10 - was empty, now becomes a populated transaction row. Give it the
# same look as the other data rows (B6:D9 style, and the right-aligned E
# style used by E6:E9/E12) by pasting formats in from those cells.
$ws.Range("B10").Value = "26.10."
$ws.Range("C10").Value = "27.10."
$ws.Range("D10").Value = "RECHNUNG VODAFONE GMBH 53408177"
$ws.Range("E10").Value = "41,03-"

$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 29.10.2024"
$ws.Range("E12").Value = "271,89-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 05.11.2024"
